$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the Binance conversion rates inside the A1 text block ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value()
$text = $text.Replace("✅ 1000 Bs = 5.38 = 21406.62 pesos", "✅ 1000 Bs = 5.35 = 21307.79 pesos")
$text = $text.Replace("✅ 21406.62 pesos = 5.34 = 945.41 Bs", "✅ 21307.79 pesos = 5.31 = 964.5 Bs")
$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 186.88
$ws2.Range("O10").Value = 3982
$ws2.Range("N12").Value = 4009.71
$ws2.Range("O12").Value = 181.5
